$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift Timestamp (col A) and Lookup (col E) forward by 3 days ---
# (24/25 Mar 2025 data window -> 27/28 Mar 2025, the next data-pull window)
$aValues = @(
    45743,45743.01041666666,45743.02083333334,45743.03125,45743.04166666666,45743.05208333334,45743.0625,45743.07291666666,45743.08333333334,45743.09375,
    45743.10416666666,45743.11458333334,45743.125,45743.13541666666,45743.14583333334,45743.15625,45743.16666666666,45743.17708333334,45743.1875,45743.19791666666,
    45743.20833333334,45743.21875,45743.22916666666,45743.23958333334,45743.25,45743.26041666666,45743.27083333334,45743.28125,45743.29166666666,45743.30208333334,
    45743.3125,45743.32291666666,45743.33333333334,45743.34375,45743.35416666666,45743.36458333334,45743.375,45743.38541666666,45743.39583333334,45743.40625,
    45743.41666666666,45743.42708333334,45743.4375,45743.44791666666,45743.45833333334,45743.46875,45743.47916666666,45743.48958333334,45743.5,45743.51041666666,
    45743.52083333334,45743.53125,45743.54166666666,45743.55208333334,45743.5625,45743.57291666666,45743.58333333334,45743.59375,45743.60416666666,45743.61458333334,
    45743.625,45743.63541666666,45743.64583333334,45743.65625,45743.66666666666,45743.67708333334,45743.6875,45743.69791666666,45743.70833333334,45743.71875,
    45743.72916666666,45743.73958333334,45743.75,45743.76041666666,45743.77083333334,45743.78125,45743.79166666666,45743.80208333334,45743.8125,45743.82291666666,
    45743.83333333334,45743.84375,45743.85416666666,45743.86458333334,45743.875,45743.88541666666,45743.89583333334,45743.90625,45743.91666666666,45743.92708333334,
    45743.9375,45743.94791666666,45743.95833333334,45743.96875,45743.97916666666,45743.98958333334,45744,45744.01041666666,45744.02083333334,45744.03125,
    45744.04166666666,45744.05208333334,45744.0625,45744.07291666666,45744.08333333334,45744.09375,45744.10416666666,45744.11458333334,45744.125,45744.13541666666,
    45744.14583333334,45744.15625,45744.16666666666,45744.17708333334,45744.1875,45744.19791666666,45744.20833333334,45744.21875,45744.22916666666,45744.23958333334,
    45744.25,45744.26041666666,45744.27083333334,45744.28125,45744.29166666666,45744.30208333334,45744.3125,45744.32291666666,45744.33333333334,45744.34375,
    45744.35416666666,45744.36458333334,45744.375,45744.38541666666,45744.39583333334,45744.40625,45744.41666666666,45744.42708333334,45744.4375,45744.44791666666,
    45744.45833333334,45744.46875,45744.47916666666,45744.48958333334,45744.5,45744.51041666666,45744.52083333334,45744.53125,45744.54166666666,45744.55208333334,
    45744.5625,45744.57291666666,45744.58333333334,45744.59375,45744.60416666666,45744.61458333334,45744.625,45744.63541666666,45744.64583333334,45744.65625,
    45744.66666666666,45744.67708333334,45744.6875,45744.69791666666,45744.70833333334,45744.71875,45744.72916666666,45744.73958333334,45744.75,45744.76041666666,
    45744.77083333334,45744.78125,45744.79166666666,45744.80208333334,45744.8125,45744.82291666666,45744.83333333334,45744.84375,45744.85416666666,45744.86458333334,
    45744.875,45744.88541666666,45744.89583333334,45744.90625,45744.91666666666,45744.92708333334,45744.9375,45744.94791666666,45744.95833333334,45744.96875,
    45744.97916666666,45744.98958333334
)

$eValues = @(
    "27.03.20251","27.03.20252","27.03.20253","27.03.20254","27.03.20255","27.03.20256","27.03.20257","27.03.20258",
    "27.03.20259","27.03.202510","27.03.202511","27.03.202512","27.03.202513","27.03.202514","27.03.202515","27.03.202516",
    "27.03.202517","27.03.202518","27.03.202519","27.03.202520","27.03.202521","27.03.202522","27.03.202523","27.03.202524",
    "27.03.202525","27.03.202526","27.03.202527","27.03.202528","27.03.202529","27.03.202530","27.03.202531","27.03.202532",
    "27.03.202533","27.03.202534","27.03.202535","27.03.202536","27.03.202537","27.03.202538","27.03.202539","27.03.202540",
    "27.03.202541","27.03.202542","27.03.202543","27.03.202544","27.03.202545","27.03.202546","27.03.202547","27.03.202548",
    "27.03.202549","27.03.202550","27.03.202551","27.03.202552","27.03.202553","27.03.202554","27.03.202555","27.03.202556",
    "27.03.202557","27.03.202558","27.03.202559","27.03.202560","27.03.202561","27.03.202562","27.03.202563","27.03.202564",
    "27.03.202565","27.03.202566","27.03.202567","27.03.202568","27.03.202569","27.03.202570","27.03.202571","27.03.202572",
    "27.03.202573","27.03.202574","27.03.202575","27.03.202576","27.03.202577","27.03.202578","27.03.202579","27.03.202580",
    "27.03.202581","27.03.202582","27.03.202583","27.03.202584","27.03.202585","27.03.202586","27.03.202587","27.03.202588",
    "27.03.202589","27.03.202590","27.03.202591","27.03.202592","27.03.202593","27.03.202594","27.03.202595","27.03.202596",
    "28.03.20251","28.03.20252","28.03.20253","28.03.20254","28.03.20255","28.03.20256","28.03.20257","28.03.20258",
    "28.03.20259","28.03.202510","28.03.202511","28.03.202512","28.03.202513","28.03.202514","28.03.202515","28.03.202516",
    "28.03.202517","28.03.202518","28.03.202519","28.03.202520","28.03.202521","28.03.202522","28.03.202523","28.03.202524",
    "28.03.202525","28.03.202526","28.03.202527","28.03.202528","28.03.202529","28.03.202530","28.03.202531","28.03.202532",
    "28.03.202533","28.03.202534","28.03.202535","28.03.202536","28.03.202537","28.03.202538","28.03.202539","28.03.202540",
    "28.03.202541","28.03.202542","28.03.202543","28.03.202544","28.03.202545","28.03.202546","28.03.202547","28.03.202548",
    "28.03.202549","28.03.202550","28.03.202551","28.03.202552","28.03.202553","28.03.202554","28.03.202555","28.03.202556",
    "28.03.202557","28.03.202558","28.03.202559","28.03.202560","28.03.202561","28.03.202562","28.03.202563","28.03.202564",
    "28.03.202565","28.03.202566","28.03.202567","28.03.202568","28.03.202569","28.03.202570","28.03.202571","28.03.202572",
    "28.03.202573","28.03.202574","28.03.202575","28.03.202576","28.03.202577","28.03.202578","28.03.202579","28.03.202580",
    "28.03.202581","28.03.202582","28.03.202583","28.03.202584","28.03.202585","28.03.202586","28.03.202587","28.03.202588",
    "28.03.202589","28.03.202590","28.03.202591","28.03.202592","28.03.202593","28.03.202594","28.03.202595","28.03.202596"
)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $aValues[$i]
    $ws.Cells.Item($r, 5).Value = $eValues[$i]
}

# --- Refreshed Notified/Actual Production (col B/C) values ---
# Rows not listed below keep their previous B/C values; "NA" = leave cell untouched.
$bcUpdates = @(
    @(22,21,"NA"),
    @(23,15,"NA"),
    @(24,23,0),
    @(25,27,5),
    @(26,106,17),
    @(27,115,25),
    @(28,122,40),
    @(29,132,61),
    @(30,236,103),
    @(31,250,120),
    @(32,292,143),
    @(33,303,173),
    @(34,428,222),
    @(35,438,278),
    @(36,452,281),
    @(37,463,305),
    @(38,568,346),
    @(39,579,407),
    @(40,589,419),
    @(41,600,439),
    @(42,656,472),
    @(43,664,511),
    @(44,671,551),
    @(45,678,659),
    @(46,705,610),
    @(47,707,608),
    @(48,708,592),
    @(49,706,636),
    @(50,667,612),
    @(51,663,497),
    @(52,658,472),
    @(53,651,477),
    @(54,594,442),
    @(55,587,436),
    @(56,575,417),
    @(57,565,431),
    @(58,476,412),
    @(59,468,361),
    @(60,457,341),
    @(61,446,351),
    @(62,311,352),
    @(63,299,307),
    @(64,288,302),
    @(65,279,277),
    @(66,145,245),
    @(67,138,186),
    @(68,129,129),
    @(69,123,87),
    @(70,10,54),
    @(71,8,22),
    @(72,7,"NA"),
    @(73,6,"NA"),
    @(76,1,"NA"),
    @(118,9,"NA"),
    @(119,10,"NA"),
    @(120,12,5),
    @(121,14,14),
    @(122,124,29),
    @(123,130,0),
    @(124,138,0),
    @(125,146,0),
    @(126,291,0),
    @(127,303,0),
    @(128,316,0),
    @(129,330,0),
    @(130,448,0),
    @(131,460,0),
    @(132,473,0),
    @(133,483,0),
    @(134,563,0),
    @(135,571,0),
    @(136,579,0),
    @(137,584,0),
    @(138,627,0),
    @(139,629,0),
    @(140,632,0),
    @(141,635,"NA"),
    @(142,631,"NA"),
    @(143,630,"NA"),
    @(144,628,"NA"),
    @(145,623,"NA"),
    @(146,579,"NA"),
    @(147,571,"NA"),
    @(148,563,"NA"),
    @(149,552,"NA"),
    @(150,492,"NA"),
    @(151,482,"NA"),
    @(152,472,"NA"),
    @(153,464,"NA"),
    @(154,370,"NA"),
    @(155,361,"NA"),
    @(156,350,"NA"),
    @(157,341,"NA"),
    @(158,235,"NA"),
    @(159,225,"NA"),
    @(160,217,"NA"),
    @(161,209,"NA"),
    @(162,114,"NA"),
    @(163,109,"NA"),
    @(164,102,"NA"),
    @(165,98,"NA"),
    @(166,10,"NA"),
    @(167,8,"NA"),
    @(168,7,"NA"),
    @(169,6,"NA"),
    @(171,0,"NA"),
    @(172,0,"NA"),
    @(173,0,"NA")
)

foreach ($u in $bcUpdates) {
    $r = $u[0]
    if ($u[1] -ne "NA") { $ws.Cells.Item($r, 2).Value = $u[1] }
    if ($u[2] -ne "NA") { $ws.Cells.Item($r, 3).Value = $u[2] }
}